# Apply the "artfynd" re-shuffle edit described by the diff.
# Rows 3, 6 and 7 exchange their species/observation data (with a couple of
# small Taxonsorteringsordning / B-column tweaks), and five other rows get a
# +1 bump to their B-column value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3  <-  old row 7 (Knärot / Goodyera repens), with B changed 99014->99015
#            and a new, empty AF3 cell.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value2 = 130937863
$ws.Range("B3").Value2 = 99015
$ws.Range("D3").Value2 = "VU"
$ws.Range("E3").Value2 = 220787
$ws.Range("F3").Value2 = "Knärot"
$ws.Range("G3").Value2 = "Goodyera repens"
$ws.Range("H3").Value2 = "(L.) R. Br."
$ws.Range("I3").Value2 = "8"
$ws.Range("J3").Value2 = "plantor/tuvor"
$ws.Range("K3").Value2 = "fullt utvecklade blad"
$ws.Range("Q3").Value2 = 489799
$ws.Range("R3").Value2 = 7004245
$ws.Range("AC3").Value2 = "Minst 8 plantor inom ca 1 m2 yta. Grävdes varsamt fram under snötäcket. Det finns sannolikt betydligt mer knärot på fyndplatsen och i skogsbeståndet där fyndplatsen ligger."
$ws.Range("AF3").Value2 = ""
$ws.Range("AH3").Value2 = "Barrskog"

# Cells that belonged to the old woodpecker record and must disappear from row 3
$ws.Range("M3").ClearContents()
$ws.Range("AJ3").ClearContents()
$ws.Range("AK3").ClearContents()
$ws.Range("AM3").ClearContents()
$ws.Range("AO3").ClearContents()

# ---------------------------------------------------------------------------
# Row 6  <-  old row 3 (Spillkråka / Dryocopus martius) content, unchanged.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value2 = 130937854
$ws.Range("B6").Value2 = 57881
$ws.Range("D6").Value2 = "NT"
$ws.Range("E6").Value2 = 100049
$ws.Range("F6").Value2 = "Spillkråka"
$ws.Range("G6").Value2 = "Dryocopus martius"
$ws.Range("H6").Value2 = "(Linnaeus, 1758)"
$ws.Range("M6").Value2 = "färska spår"
$ws.Range("Q6").Value2 = 489668
$ws.Range("R6").Value2 = 7004128
$ws.Range("AC6").Value2 = "Rejäla hackspår, färska och äldre, I två levande granar och i ytlig grov rotdel."
$ws.Range("AH6").Value2 = "Granskog"
$ws.Range("AJ6").Value2 = "gran"
$ws.Range("AK6").Value2 = "Picea abies"
$ws.Range("AM6").Value2 = "Trädstam på levande träd"
$ws.Range("AO6").Value2 = "Stem on living tree # Picea abies"

# Cells that belonged to the old Revlummer record and must disappear from row 6
$ws.Range("J6").ClearContents()
$ws.Range("AF6").ClearContents()

# ---------------------------------------------------------------------------
# Row 7  <-  old row 6 (Revlummer / Lycopodium annotinum), with B changed
#            97879->97880.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value2 = 130937857
$ws.Range("B7").Value2 = 97880
$ws.Range("D7").Value2 = "LC"
$ws.Range("E7").Value2 = 221945
$ws.Range("F7").Value2 = "Revlummer"
$ws.Range("G7").Value2 = "Lycopodium annotinum"
$ws.Range("H7").Value2 = "L."
$ws.Range("Q7").Value2 = 489680
$ws.Range("R7").Value2 = 7004154
$ws.Range("AH7").Value2 = "Barrskog"

# Cell that belonged to the old Knärot record and must disappear from row 7
$ws.Range("AC7").ClearContents()

# ---------------------------------------------------------------------------
# Small Taxonsorteringsordning (column B) bumps on five other rows.
# ---------------------------------------------------------------------------
$ws.Range("B11").Value2 = 97880
$ws.Range("B12").Value2 = 99352
$ws.Range("B17").Value2 = 97880
$ws.Range("B20").Value2 = 97880
$ws.Range("B21").Value2 = 97880
